# "first version of streamlit page"
# On the "project_returns" sheet, the return value for project 1120 is
# corrected from 12000 to 120000, and the active selection is left on C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("project_returns")
$ws.Activate()

$ws.Range("C2").Value = 120000

$ws.Range("C3").Select()
